# Updates odds values in the "Jogos da Semana" worksheet to reflect the
# latest FlashScore odds refresh for 2024-10-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Corinthians vs Athletico-PR)
$ws.Range("O2").Value = 1.33
$ws.Range("P2").Value = 3.4

# Row 3 (Flamengo RJ vs Fluminense)
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 2.75
$ws.Range("N3").Value = 7.5
$ws.Range("X3").Value = 8.5
$ws.Range("AH3").Value = 9
$ws.Range("AJ3").Value = 13
$ws.Range("AL3").Value = 34
$ws.Range("AP3").Value = 26
$ws.Range("AW3").Value = 5.5

# Row 10 (Sport Huancayo vs Grau)
$ws.Range("Q10").Value = 2.1
$ws.Range("R10").Value = 1.7

# Row 13 (Nacional vs Miramar)
$ws.Range("M13").Value = 1.05
$ws.Range("N13").Value = 11
$ws.Range("Q13").Value = 1.9
$ws.Range("R13").Value = 1.95
